$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 699.8
$ws.Range("I2").Value = 750
$ws.Range("K2").Value = 750
$ws.Range("M2").Value = -637
$ws.Range("H12").Value = 147
$ws.Range("I12").Value = 129.8
$ws.Range("K12").Value = 129.8
$ws.Range("M12").Value = 40.19999999999999
$ws.Range("H18").Value = 333.2857
$ws.Range("I18").Value = 333.2857
$ws.Range("K18").Value = 333.2857
$ws.Range("M18").Value = -49.28570000000002
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H29").Value = 4600.231
$ws.Range("J29").Value = 6388.8887
$ws.Range("L29").Value = 19166.6661
$ws.Range("N29").Value = -19728.6661
$ws.Range("H31").Value = 282.7143
$ws.Range("I31").Value = 282.7143
$ws.Range("K31").Value = 848.1428999999999
$ws.Range("M31").Value = -618.1428999999999
$ws.Range("H40").Value = 5080
$ws.Range("I40").Value = 4100
$ws.Range("K40").Value = 4100
$ws.Range("M40").Value = -3925
$ws.Range("H96").Value = 1739.5454
$ws.Range("I96").Value = 1123.5714
$ws.Range("J96").Value = 2817.5
$ws.Range("K96").Value = 3370.7142
$ws.Range("L96").Value = 8452.5
$ws.Range("M96").Value = -1997.7142
$ws.Range("N96").Value = -11198.5
$ws.Range("H136").Value = 199999
$ws.Range("J136").Value = 199999
$ws.Range("L136").Value = 199999
$ws.Range("N136").Value = -210199
$ws.Range("H137").Value = 6408.593
$ws.Range("I137").Value = 1982.1
$ws.Range("J137").Value = 19055.715
$ws.Range("K137").Value = 5946.299999999999
$ws.Range("L137").Value = 57167.145
$ws.Range("M137").Value = -3396.299999999999
$ws.Range("N137").Value = -62267.145
$ws.Range("H138").Value = 6838.5
$ws.Range("I138").Value = 1371.5
$ws.Range("J138").Value = 9962.5
$ws.Range("K138").Value = 4114.5
$ws.Range("L138").Value = 29887.5
$ws.Range("M138").Value = 1025.5
$ws.Range("N138").Value = -40167.5
$ws.Range("H139").Value = 120000
$ws.Range("J139").Value = 120000
$ws.Range("L139").Value = 120000
$ws.Range("N139").Value = -130280

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2298.1724
$ws.Range("I2").Value = 1730.8334
$ws.Range("K2").Value = 1730.8334
$ws.Range("M2").Value = -1617.8334
$ws.Range("H32").Value = 10870979
$ws.Range("I32").Value = 11765932
$ws.Range("J32").Value = 3689.8572
$ws.Range("K32").Value = 11765932
$ws.Range("L32").Value = 3689.8572
$ws.Range("M32").Value = -11765645
$ws.Range("N32").Value = -4263.8572
$ws.Range("H45").Value = 4138.8
$ws.Range("I45").Value = 3777.8
$ws.Range("J45").Value = 4499.8
$ws.Range("K45").Value = 3777.8
$ws.Range("L45").Value = 4499.8
$ws.Range("M45").Value = -3400.8
$ws.Range("N45").Value = -5253.8
$ws.Range("H116").Value = 2298.1724
$ws.Range("I116").Value = 1730.8334
$ws.Range("K116").Value = 1730.8334
$ws.Range("M116").Value = 563.1666
$ws.Range("H132").Value = 41726576
$ws.Range("I132").Value = 1902.6
$ws.Range("K132").Value = 5707.799999999999
$ws.Range("M132").Value = -3177.799999999999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2298.1724
$ws.Range("I3").Value = 1730.8334
$ws.Range("K3").Value = 1730.8334
$ws.Range("M3").Value = -1616.8334
$ws.Range("H20").Value = 1968.175
$ws.Range("I20").Value = 1696.75
$ws.Range("J20").Value = 2375.3125
$ws.Range("K20").Value = 1696.75
$ws.Range("L20").Value = 2375.3125
$ws.Range("M20").Value = -1449.75
$ws.Range("N20").Value = -2869.3125
$ws.Range("H22").Value = 299.92856
$ws.Range("I22").Value = 299.92856
$ws.Range("K22").Value = 299.92856
$ws.Range("M22").Value = -126.92856
$ws.Range("H57").Value = 199999
$ws.Range("J57").Value = 199999
$ws.Range("L57").Value = 199999
$ws.Range("N57").Value = -201439
$ws.Range("H107").Value = 2572.2273
$ws.Range("I107").Value = 1625.8125
$ws.Range("K107").Value = 1625.8125
$ws.Range("M107").Value = 294.1875
$ws.Range("H133").Value = 60000
$ws.Range("I133").Value = 60000
$ws.Range("K133").Value = 60000
$ws.Range("M133").Value = -54940
$ws.Range("H136").Value = 199999
$ws.Range("J136").Value = 199999
$ws.Range("L136").Value = 199999
$ws.Range("N136").Value = -210199

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1216.5555
$ws.Range("I16").Value = 1278.4286
$ws.Range("K16").Value = 1278.4286
$ws.Range("M16").Value = -991.4286
$ws.Range("H59").Value = 50115
$ws.Range("J59").Value = 50115
$ws.Range("L59").Value = 50115
$ws.Range("N59").Value = -52405
$ws.Range("H62").Value = 3150.9473
$ws.Range("I62").Value = 3177.389
$ws.Range("J62").Value = 2675
$ws.Range("K62").Value = 3177.389
$ws.Range("L62").Value = 2675
$ws.Range("M62").Value = -2553.389
$ws.Range("N62").Value = -3923
$ws.Range("H65").Value = 3150.9473
$ws.Range("I65").Value = 3177.389
$ws.Range("J65").Value = 2675
$ws.Range("K65").Value = 15886.945
$ws.Range("L65").Value = 13375
$ws.Range("M65").Value = -12766.945
$ws.Range("N65").Value = -19615
$ws.Range("H113").Value = 1216.5555
$ws.Range("I113").Value = 1278.4286
$ws.Range("K113").Value = 1278.4286
$ws.Range("M113").Value = 891.5714
$ws.Range("H132").Value = 2412.5833
$ws.Range("I132").Value = 2457.1042
$ws.Range("J132").Value = 2234.5
$ws.Range("K132").Value = 7371.312600000001
$ws.Range("L132").Value = 6703.5
$ws.Range("M132").Value = -4841.312600000001
$ws.Range("N132").Value = -11763.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 1000
$ws.Range("I115").Value = 1000
$ws.Range("K115").Value = 3000
$ws.Range("M115").Value = -1825
$ws.Range("H128").Value = 187315
$ws.Range("I128").Value = 187315
$ws.Range("K128").Value = 561945
$ws.Range("M128").Value = -556965
$ws.Range("H132").Value = 4765994
$ws.Range("I132").Value = 1570.8572
$ws.Range("J132").Value = 9530417
$ws.Range("K132").Value = 14137.7148
$ws.Range("L132").Value = 85773753
$ws.Range("M132").Value = -11607.7148
$ws.Range("N132").Value = -85778813

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 199999
$ws.Range("J140").Value = 199999
$ws.Range("L140").Value = 199999
$ws.Range("N140").Value = -210359

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2403.5334
$ws.Range("J22").Value = 2543.1428
$ws.Range("L22").Value = 2543.1428
$ws.Range("N22").Value = -3133.1428
$ws.Range("H27").Value = 2403.5334
$ws.Range("J27").Value = 2543.1428
$ws.Range("L27").Value = 2543.1428
$ws.Range("N27").Value = -2757.1428
$ws.Range("H46").Value = 2218.2258
$ws.Range("H61").Value = 6806.533
$ws.Range("I61").Value = 6255.5
$ws.Range("J61").Value = 7908.6
$ws.Range("K61").Value = 6255.5
$ws.Range("L61").Value = 7908.6
$ws.Range("M61").Value = -6053.5
$ws.Range("N61").Value = -8312.6
$ws.Range("H113").Value = 6806.533
$ws.Range("I113").Value = 6255.5
$ws.Range("J113").Value = 7908.6
$ws.Range("K113").Value = 6255.5
$ws.Range("L113").Value = 7908.6
$ws.Range("M113").Value = -4085.5
$ws.Range("N113").Value = -12248.6
$ws.Range("H122").Value = 5004471
$ws.Range("J122").Value = 13893922
$ws.Range("L122").Value = 41681766
$ws.Range("N122").Value = -41686666

Write-Host "Applied all cell updates"